# NPX09 initial commit (Reformed list, dilution heatmaps)
# Rows 81-86 in the "Elution" column (F) were re-mapped from "FPLC" to
# "None", and the sheet's active selection moved down to F87.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F81:F86 -> "None" (was "FPLC")
$ws.Range("F81:F86").Value = "None"

# Move the live selection to F87 (single cell), matching where the author
# left the cursor after editing this block of rows.
$ws.Range("F87").Select() | Out-Null
